# Generate Report for Handoff
#
# The localization-status report is regenerated: rows in each sheet are
# re-ordered by latest handoff date, and the row for
# "311bc6e3-95de-4af4-ae66-119201a72867.md" moves from "Handed back: in
# sync with en-US" to "Ready for handoff" with refreshed handoff
# file/datetime info. Hyperlink targets stay anchored to their original
# cell position (same r:id), only the visible text changes - so we update
# cell values *and* the matching hyperlink's displayed text in lock-step.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkMap($ws) {
    $map = @{}
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        $map[$addr] = $h
    }
    return $map
}

function AddrKey($addr) {
    if ($addr -match '^([A-Z]+)([0-9]+)$') {
        return '$' + $matches[1] + '$' + $matches[2]
    }
    return $addr
}

function Set-CellText($ws, $map, $addr, $text) {
    $ws.Range($addr).Value = $text
    $key = AddrKey $addr
    if ($map.ContainsKey($key)) {
        $map[$key].TextToDisplay = $text
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$map1 = Get-HyperlinkMap $ws1

Set-CellText $ws1 $map1 "A2" "ffffa1923a7b-fa80-43ac-9549-8f327d1787e2.md"
Set-CellText $ws1 $map1 "B2" "Handed back: in sync with en-US"
Set-CellText $ws1 $map1 "C2" "Handed back: in sync with en-US"
Set-CellText $ws1 $map1 "D2" "2016-03-22 07:08:57"

Set-CellText $ws1 $map1 "A3" "ffffff6a00001f-4ff7-4398-ab3c-e50ae7d837ce.md"
Set-CellText $ws1 $map1 "B3" "Handed back: in sync with en-US"
Set-CellText $ws1 $map1 "C3" "Handed back: in sync with en-US"
Set-CellText $ws1 $map1 "D3" "2016-03-22 07:08:57"

Set-CellText $ws1 $map1 "A4" "311bc6e3-95de-4af4-ae66-119201a72867.md"
Set-CellText $ws1 $map1 "B4" "Ready for handoff"
Set-CellText $ws1 $map1 "C4" "Ready for handoff"
Set-CellText $ws1 $map1 "D4" "2016-03-22 07:13:01"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$map2 = Get-HyperlinkMap $ws2

Set-CellText $ws2 $map2 "A2" "ffffa1923a7b-fa80-43ac-9549-8f327d1787e2.md"
Set-CellText $ws2 $map2 "C2" "Handed back: in sync with en-US"
Set-CellText $ws2 $map2 "D2" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"
Set-CellText $ws2 $map2 "E2" "2016-03-22 07:08:53"
Set-CellText $ws2 $map2 "F2" "9993c348-c562-422b-8d38-0d8a9c505173.md"
Set-CellText $ws2 $map2 "G2" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"
Set-CellText $ws2 $map2 "H2" "2016-03-22 07:09:15"

Set-CellText $ws2 $map2 "A3" "ffffff6a00001f-4ff7-4398-ab3c-e50ae7d837ce.md"
Set-CellText $ws2 $map2 "C3" "Handed back: in sync with en-US"
Set-CellText $ws2 $map2 "D3" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"
Set-CellText $ws2 $map2 "E3" "2016-03-22 07:08:53"
Set-CellText $ws2 $map2 "F3" "9993c348-c562-422b-8d38-0d8a9c505173.md"
Set-CellText $ws2 $map2 "G3" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"
Set-CellText $ws2 $map2 "H3" "2016-03-22 07:09:15"

Set-CellText $ws2 $map2 "A4" "311bc6e3-95de-4af4-ae66-119201a72867.md"
Set-CellText $ws2 $map2 "C4" "Ready for handoff"
Set-CellText $ws2 $map2 "D4" "311bc6e3-95de-4af4-ae66-119201a72867.07072297964b3311d1a11a79a5f2a453d2cd8058.zh-cn.xlf"
Set-CellText $ws2 $map2 "E4" "2016-03-22 07:12:57"
Set-CellText $ws2 $map2 "F4" "311bc6e3-95de-4af4-ae66-119201a72867.md"
Set-CellText $ws2 $map2 "G4" "311bc6e3-95de-4af4-ae66-119201a72867.07072297964b3311d1a11a79a5f2a453d2cd8058.zh-cn.xlf"
Set-CellText $ws2 $map2 "H4" "2016-03-22 07:12:24"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$map3 = Get-HyperlinkMap $ws3

Set-CellText $ws3 $map3 "A2" "ffffa1923a7b-fa80-43ac-9549-8f327d1787e2.md"
Set-CellText $ws3 $map3 "C2" "Handed back: in sync with en-US"
Set-CellText $ws3 $map3 "D2" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"
Set-CellText $ws3 $map3 "E2" "2016-03-22 07:08:57"
Set-CellText $ws3 $map3 "F2" "9993c348-c562-422b-8d38-0d8a9c505173.md"
Set-CellText $ws3 $map3 "G2" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"
Set-CellText $ws3 $map3 "H2" "2016-03-22 07:09:21"

Set-CellText $ws3 $map3 "A3" "ffffff6a00001f-4ff7-4398-ab3c-e50ae7d837ce.md"
Set-CellText $ws3 $map3 "C3" "Handed back: in sync with en-US"
Set-CellText $ws3 $map3 "D3" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"
Set-CellText $ws3 $map3 "E3" "2016-03-22 07:08:57"
Set-CellText $ws3 $map3 "F3" "9993c348-c562-422b-8d38-0d8a9c505173.md"
Set-CellText $ws3 $map3 "G3" "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"
Set-CellText $ws3 $map3 "H3" "2016-03-22 07:09:21"

Set-CellText $ws3 $map3 "A4" "311bc6e3-95de-4af4-ae66-119201a72867.md"
Set-CellText $ws3 $map3 "C4" "Ready for handoff"
Set-CellText $ws3 $map3 "D4" "311bc6e3-95de-4af4-ae66-119201a72867.07072297964b3311d1a11a79a5f2a453d2cd8058.de-de.xlf"
Set-CellText $ws3 $map3 "E4" "2016-03-22 07:13:01"
Set-CellText $ws3 $map3 "F4" "311bc6e3-95de-4af4-ae66-119201a72867.md"
Set-CellText $ws3 $map3 "G4" "311bc6e3-95de-4af4-ae66-119201a72867.07072297964b3311d1a11a79a5f2a453d2cd8058.de-de.xlf"
Set-CellText $ws3 $map3 "H4" "2016-03-22 07:12:30"
